$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 444, shifting existing rows 444:473 down to 445:474
$ws.Rows.Item(444).Insert()

# Populate the newly inserted row 444 with the new record
$ws.Cells.Item(444, 1).Value = 11
$ws.Cells.Item(444, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(444, 3).Value = "Bíobío"
$ws.Cells.Item(444, 4).Value = "2023-12-07"
$ws.Cells.Item(444, 5).Value = 8
$ws.Cells.Item(444, 6).Value = 100112045
$ws.Cells.Item(444, 7).Value = "Zapallo"
$ws.Cells.Item(444, 8).Value = "Paine"
$ws.Cells.Item(444, 9).Value = "1a nueva(o)"
$ws.Cells.Item(444, 10).Value = 400
$ws.Cells.Item(444, 11).Value = 1000
$ws.Cells.Item(444, 12).Value = 1000
$ws.Cells.Item(444, 13).Value = 1000
$ws.Cells.Item(444, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(444, 15).Value = "Región Metropolitana"
$ws.Cells.Item(444, 16).Value = 1000
$ws.Cells.Item(444, 17).Value = 1
$ws.Cells.Item(444, 18).Value = "Hortaliza"
